$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.272.93"
$ws.Range("E2").Value = "  -0.78%  "

# Row 3
$ws.Range("D3").Value = "1.885.33"
$ws.Range("E3").Value = "  -1.42%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4672"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.72%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2829"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.50%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "

# Row 10
$ws.Range("E10").Value = "  +4.99%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07772"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.03%  "

# Row 13
$ws.Range("D13").Value = "1.888.44"
$ws.Range("E13").Value = "  -1.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.122"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6682"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.28%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "283.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +10.50%  "

# Row 17
$ws.Range("D17").Value = "30.279.82"
$ws.Range("E17").Value = "  -0.83%  "

# Row 18
$ws.Range("E18").Value = "  +0.05%  "

# Row 19
$ws.Range("E19").Value = "  -0.49%  "

# Row 20
$ws.Range("D20").Value = "2.133.00"
$ws.Range("E20").Value = "  -1.32%  "

# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007304"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.42%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.359"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.78%  "

# Row 23
$ws.Range("E23").Value = "  +0.10%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.171"
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.343"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "

# Row 28
$ws.Range("E28").Value = "  -3.16%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.373"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09733"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.38%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.457"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.52%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.483"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.172"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04702"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.54%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7084"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.75%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.094"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.58%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.712"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01867"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.56%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.646"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.67%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.526"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.40%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.44%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.973"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.98%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8694"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.97%  "

# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.11%  "

# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4200"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.99%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "987.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.215"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.246"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.59%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1163"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.94%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.12%  "
